$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (matches original inlineStr text cells)
function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "246.27"
Set-TextCell $ws "E2" "0.23%"
Set-TextCell $ws "D3" "26.05"
Set-TextCell $ws "E3" "2.65%"
Set-TextCell $ws "D4" "5.193"
Set-TextCell $ws "E4" "2.96%"
Set-TextCell $ws "D5" "0.05591"
Set-TextCell $ws "D6" "6.481"
Set-TextCell $ws "E6" "-1.35%"
Set-TextCell $ws "D7" "0.8127"
Set-TextCell $ws "E7" "-0.37%"
Set-TextCell $ws "D8" "0.8445"
Set-TextCell $ws "E8" "1.14%"
Set-TextCell $ws "D9" "0.06939"
Set-TextCell $ws "E9" "-0.16%"
Set-TextCell $ws "D10" "0.02837"
Set-TextCell $ws "E10" "0.42%"
Set-TextCell $ws "D11" "0.09380"
Set-TextCell $ws "E11" "-0.24%"
Set-TextCell $ws "D12" "0.001516"
Set-TextCell $ws "E12" "-0.29%"
$ws.Range("B13").Value = "One"
$ws.Range("C13").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell $ws "D13" "0.0005999"
Set-TextCell $ws "E13" "0.17%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws "D14" "0.006117"
Set-TextCell $ws "E14" "-0.98%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D15" "3.606"
Set-TextCell $ws "E15" "3.12%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell $ws "D16" "3.028"
Set-TextCell $ws "E16" "0.62%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell $ws "D17" "2.058"
Set-TextCell $ws "E17" "-1.58%"
Set-TextCell $ws "D19" "0.1328"
Set-TextCell $ws "E19" "-0.73%"
Set-TextCell $ws "D20" "0.03123"
Set-TextCell $ws "E20" "-3.47%"
Set-TextCell $ws "D21" "0.1294"
Set-TextCell $ws "E21" "-1.97%"
Set-TextCell $ws "D22" "3.768"
Set-TextCell $ws "E22" "0.82%"
Set-TextCell $ws "D23" "0.04653"
Set-TextCell $ws "E23" "-0.77%"
Set-TextCell $ws "D24" "0.1374"
Set-TextCell $ws "E24" "2.47%"
Set-TextCell $ws "E25" "0.10%"
Set-TextCell $ws "D26" "0.004546"
Set-TextCell $ws "E26" "5.94%"
Set-TextCell $ws "D27" "0.00009598"
Set-TextCell $ws "E27" "-1.02%"
Set-TextCell $ws "E28" "-27.55%"
Set-TextCell $ws "D40" "0.03651"
Set-TextCell $ws "E40" "-0.29%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell $ws "D41" "0.1367"
Set-TextCell $ws "E41" "0.32%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell $ws "D42" "0.002649"
Set-TextCell $ws "E42" "-3.63%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell $ws "D43" "0.003411"
Set-TextCell $ws "E43" "-44.88%"
Set-TextCell $ws "D44" "0.007980"
Set-TextCell $ws "E44" "7.40%"
Set-TextCell $ws "D45" "0.00005385"
Set-TextCell $ws "E45" "1.73%"
Set-TextCell $ws "E46" "0.00%"
Set-TextCell $ws "D47" "0.1450"
Set-TextCell $ws "E47" "-19.44%"
Set-TextCell $ws "D48" "0.002410"
Set-TextCell $ws "E48" "19.57%"
Set-TextCell $ws "D49" "0.00002100"
Set-TextCell $ws "E49" "0.00%"
Set-TextCell $ws "D50" "0.0002000"
Set-TextCell $ws "E50" "0.00%"
